$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the degree distribution values for rows 2-5
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

# Remove the now-invalid extra rows (6-19) that belonged to the old distribution
$ws.Range("A6:B19").ClearContents()
